$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44774
$ws.Range("E2").Value = 2253
$ws.Range("F2").Value = 2253
$ws.Range("G2").Value = 1217
$ws.Range("H2").Value = 833
$ws.Range("I2").Value = 688
$ws.Range("J2").Value = 145
$ws.Range("K2").Value = 59365
$ws.Range("L2").Value = 36636
$ws.Range("M2").Value = 22729
$ws.Range("N2").Value = 21806
$ws.Range("O2").Value = 923
$ws.Range("P2").Value = 3769
$ws.Range("Q2").Value = 5276
$ws.Range("R2").Value = -1652
$ws.Range("S2").Value = -4037
$ws.Range("T2").Value = 204
$ws.Range("U2").Value = 5072
$ws.Range("V2").Value = 19561
$ws.Range("W2").Value = 5.03
$ws.Range("X2").Value = 1.86
$ws.Range("Y2").Value = 3.2
$ws.Range("Z2").Value = 1.35
$ws.Range("AA2").Value = 161.18
$ws.Range("AB2").Value = 505.41
$ws.Range("AC2").Value = 913
$ws.Range("AD2").Value = 44.43
$ws.Range("AE2").Value = 29634
$ws.Range("AF2").Value = 1.37
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = 32.07
$ws.Range("AJ2").Value = 75384180

$ws.Range("D3").Value = 46026
$ws.Range("E3").Value = 3895
$ws.Range("F3").Value = 3895
$ws.Range("G3").Value = 3243
$ws.Range("H3").Value = 2386
$ws.Range("I3").Value = 2168
$ws.Range("J3").Value = 218
$ws.Range("K3").Value = 55101
$ws.Range("L3").Value = 30391
$ws.Range("M3").Value = 24710
$ws.Range("N3").Value = 23576
$ws.Range("O3").Value = 1134
$ws.Range("P3").Value = 3769
$ws.Range("Q3").Value = 10284
$ws.Range("R3").Value = -643
$ws.Range("S3").Value = -6645
$ws.Range("T3").Value = 239
$ws.Range("U3").Value = 10045
$ws.Range("V3").Value = 12989
$ws.Range("W3").Value = 8.46
$ws.Range("X3").Value = 5.18
$ws.Range("Y3").Value = 9.55
$ws.Range("Z3").Value = 4.17
$ws.Range("AA3").Value = 122.99
$ws.Range("AB3").Value = 556.11
$ws.Range("AC3").Value = 2876
$ws.Range("AD3").Value = 14.15
$ws.Range("AE3").Value = 32039
$ws.Range("AF3").Value = 1.27
$ws.Range("AG3").Value = 500
$ws.Range("AH3").Value = 1.23
$ws.Range("AI3").Value = 16.97
$ws.Range("AJ3").Value = 75384180

$ws.Range("D4").Value = 47499
$ws.Range("E4").Value = 5172
$ws.Range("F4").Value = 5172
$ws.Range("G4").Value = 4595
$ws.Range("H4").Value = 3310
$ws.Range("I4").Value = 3067
$ws.Range("J4").Value = 242
$ws.Range("K4").Value = 57846
$ws.Range("L4").Value = 30297
$ws.Range("M4").Value = 27549
$ws.Range("N4").Value = 26220
$ws.Range("O4").Value = 1329
$ws.Range("P4").Value = 3769
$ws.Range("Q4").Value = 8398
$ws.Range("R4").Value = -890
$ws.Range("S4").Value = -1458
$ws.Range("T4").Value = 498
$ws.Range("U4").Value = 7899
$ws.Range("V4").Value = 11987
$ws.Range("W4").Value = 10.89
$ws.Range("X4").Value = 6.97
$ws.Range("Y4").Value = 12.32
$ws.Range("Z4").Value = 5.86
$ws.Range("AA4").Value = 109.97
$ws.Range("AB4").Value = 627.3
$ws.Range("AC4").Value = 4069
$ws.Range("AD4").Value = 11.58
$ws.Range("AE4").Value = 35632
$ws.Range("AF4").Value = 1.32
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 1.49
$ws.Range("AI4").Value = 16.79
$ws.Range("AJ4").Value = 75384180

$ws.Range("D5").Value = 15210
$ws.Range("E5").Value = 856
$ws.Range("F5").Value = 856
$ws.Range("G5").Value = 484
$ws.Range("H5").Value = 4137
$ws.Range("I5").Value = 3896
$ws.Range("J5").Value = 240
$ws.Range("K5").Value = 65394
$ws.Range("L5").Value = 35770
$ws.Range("M5").Value = 29625
$ws.Range("N5").Value = 28034
$ws.Range("O5").Value = 1591
$ws.Range("P5").Value = 3769
$ws.Range("Q5").Value = 4532
$ws.Range("R5").Value = -2290
$ws.Range("S5").Value = -976
$ws.Range("T5").Value = 2112
$ws.Range("U5").Value = 2421
$ws.Range("V5").Value = 13169
$ws.Range("W5").Value = 5.63
$ws.Range("X5").Value = 27.2
$ws.Range("Y5").Value = 14.36
$ws.Range("Z5").Value = 6.71
$ws.Range("AA5").Value = 120.74
$ws.Range("AB5").Value = 717.2
$ws.Range("AC5").Value = 5169
$ws.Range("AD5").Value = 7.82
$ws.Range("AE5").Value = 40000
$ws.Range("AF5").Value = 1.01
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 2.47
$ws.Range("AI5").Value = 17.99
$ws.Range("AJ5").Value = 75384180

$ws.Range("D6").Value = 15478
$ws.Range("E6").Value = 954
$ws.Range("F6").Value = 954
$ws.Range("G6").Value = 1141
$ws.Range("H6").Value = 9313
$ws.Range("I6").Value = 9171
$ws.Range("K6").Value = 34097
$ws.Range("L6").Value = 12963
$ws.Range("M6").Value = 21134
$ws.Range("N6").Value = 19217
$ws.Range("P6").Value = 2987
$ws.Range("Q6").Value = 1279
$ws.Range("R6").Value = -1322
$ws.Range("S6").Value = -9985
$ws.Range("T6").Value = 964
$ws.Range("U6").Value = 316
$ws.Range("V6").Value = 7650
$ws.Range("W6").Value = 6.17
$ws.Range("X6").Value = 60.17
$ws.Range("Y6").Value = 38.82
$ws.Range("Z6").Value = 18.72
$ws.Range("AA6").Value = 61.34
$ws.Range("AB6").Value = 1268.89
$ws.Range("AC6").Value = 16984
$ws.Range("AD6").Value = 1.02
$ws.Range("AE6").Value = 33408
$ws.Range("AF6").Value = 0.52
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 0.87
$ws.Range("AI6").Value = 0.94
$ws.Range("AJ6").Value = 59741721

$ws.Range("D7").Value = 16348
$ws.Range("E7").Value = 1378
$ws.Range("G7").Value = 3086
$ws.Range("H7").Value = 2729
$ws.Range("I7").Value = 2519
$ws.Range("K7").Value = 38196
$ws.Range("L7").Value = 14600
$ws.Range("M7").Value = 23597
$ws.Range("N7").Value = 21524
$ws.Range("P7").Value = 2987
$ws.Range("Q7").Value = -456
$ws.Range("R7").Value = -1774
$ws.Range("S7").Value = 106
$ws.Range("T7").Value = 2014
$ws.Range("U7").Value = -2898
$ws.Range("W7").Value = 8.43
$ws.Range("X7").Value = 16.69
$ws.Range("Y7").Value = 12.37
$ws.Range("Z7").Value = 7.55
$ws.Range("AA7").Value = 61.87
$ws.Range("AC7").Value = 4216
$ws.Range("AD7").Value = 2.38
$ws.Range("AE7").Value = 37418
$ws.Range("AF7").Value = 0.27
$ws.Range("AG7").Value = 225
$ws.Range("AH7").Value = 2.24
$ws.Range("AI7").Value = 5.34

$ws.Range("D8").Value = 17007
$ws.Range("E8").Value = 1524
$ws.Range("G8").Value = 2732
$ws.Range("H8").Value = 2126
$ws.Range("I8").Value = 1992
$ws.Range("K8").Value = 39472
$ws.Range("L8").Value = 13958
$ws.Range("M8").Value = 25514
$ws.Range("N8").Value = 23312
$ws.Range("P8").Value = 2987
$ws.Range("Q8").Value = 1624
$ws.Range("R8").Value = -1231
$ws.Range("S8").Value = -1128
$ws.Range("T8").Value = 531
$ws.Range("U8").Value = 291
$ws.Range("W8").Value = 8.96
$ws.Range("X8").Value = 12.5
$ws.Range("Y8").Value = 8.88
$ws.Range("Z8").Value = 5.47
$ws.Range("AA8").Value = 54.71
$ws.Range("AC8").Value = 3334
$ws.Range("AD8").Value = 3.01
$ws.Range("AE8").Value = 40526
$ws.Range("AF8").Value = 0.25
$ws.Range("AG8").Value = 275
$ws.Range("AH8").Value = 2.74
$ws.Range("AI8").Value = 8.25

$ws.Range("D9").Value = 17316
$ws.Range("E9").Value = 1596
$ws.Range("G9").Value = 2904
$ws.Range("H9").Value = 2238
$ws.Range("I9").Value = 2088
$ws.Range("K9").Value = 41802
$ws.Range("L9").Value = 14286
$ws.Range("M9").Value = 27516
$ws.Range("N9").Value = 25172
$ws.Range("P9").Value = 2987
$ws.Range("Q9").Value = 2970
$ws.Range("R9").Value = -816
$ws.Range("S9").Value = -1346
$ws.Range("T9").Value = 574
$ws.Range("U9").Value = -701
$ws.Range("W9").Value = 9.22
$ws.Range("X9").Value = 12.93
$ws.Range("Y9").Value = 8.61
$ws.Range("Z9").Value = 5.51
$ws.Range("AA9").Value = 51.92
$ws.Range("AC9").Value = 3495
$ws.Range("AD9").Value = 2.88
$ws.Range("AE9").Value = 43760
$ws.Range("AF9").Value = 0.23
$ws.Range("AG9").Value = 275
$ws.Range("AH9").Value = 2.74
$ws.Range("AI9").Value = 7.87
